$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C) for rows 2-11
# from serial date 45172 (2023-09-03) to 45175 (2023-09-06)
for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45172) {
        $cell.Value2 = 45175
    }
}
